$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Agriculture & food"
$ws.Range("B3").Value = "Petrochemicals"
$ws.Range("B7").Value = "Agriculture & food"
$ws.Range("B9").Value = "Petrochemicals"

$ws.Range("B10").Select()
